$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Re-style the three balance-sheet tables (slides 14-16) with the new
#    built-in table style that replaced the old default table style.
# ---------------------------------------------------------------------------
$newTableStyleId = "{1DE81776-6C28-4B40-B7C1-E9970BEC1806}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($i = 1; $i -le $s.Shapes.Count; $i++) {
        $sh = $s.Shapes.Item($i)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle($newTableStyleId)
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Swap the presentation's design theme away from "Integral" / Red Violet
#    back to the default "Office Theme" colours (the deck's other theme
#    part, used only by the notes master, already carries the Office
#    palette, so the live design theme now takes on those RGB values).
# ---------------------------------------------------------------------------
function RGBInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2),16)
    $g = [Convert]::ToInt32($hex.Substring(2,2),16)
    $b = [Convert]::ToInt32($hex.Substring(4,2),16)
    return $b * 65536 + $g * 256 + $r
}

# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (theme colour order 1-12)
$officeTheme = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")

$slide1 = $p.Slides.Item(1)
$tcs = $slide1.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = RGBInt($officeTheme[$i - 1])
}
